$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4862.387
$ws.Range("I132").Value = 3852.8372
$ws.Range("J132").Value = 7147.1577
$ws.Range("K132").Value = 11558.5116
$ws.Range("L132").Value = 21441.4731
$ws.Range("M132").Value = -9028.5116
$ws.Range("N132").Value = -26501.4731
$ws.Range("H137").Value = 1993.95
$ws.Range("I137").Value = 2041.0769
$ws.Range("K137").Value = 6123.2307
$ws.Range("M137").Value = -3573.2307
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6403.8213
$ws.Range("I32").Value = 6522.806
$ws.Range("K32").Value = 6522.806
$ws.Range("M32").Value = -6235.806
$ws.Range("H132").Value = 3728.6416
$ws.Range("I132").Value = 1431.6
$ws.Range("J132").Value = 8195.111000000001
$ws.Range("K132").Value = 4294.799999999999
$ws.Range("L132").Value = 24585.333
$ws.Range("M132").Value = -1764.799999999999
$ws.Range("N132").Value = -29645.333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 400
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 700
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 1220
$ws.Range("N107").Value = -3940
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 38000
$ws.Range("J135").Value = 38000
$ws.Range("L135").Value = 38000
$ws.Range("N135").Value = -48140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4911.1
$ws.Range("I16").Value = 5302.2
$ws.Range("J16").Value = 4520
$ws.Range("K16").Value = 5302.2
$ws.Range("L16").Value = 4520
$ws.Range("M16").Value = -5015.2
$ws.Range("N16").Value = -5094
$ws.Range("H31").Value = 6062559
$ws.Range("I31").Value = 1470.8223
$ws.Range("J31").Value = 33337456
$ws.Range("K31").Value = 1470.8223
$ws.Range("L31").Value = 33337456
$ws.Range("M31").Value = -1175.8223
$ws.Range("N31").Value = -33338046
$ws.Range("H34").Value = 6062559
$ws.Range("I34").Value = 1470.8223
$ws.Range("J34").Value = 33337456
$ws.Range("K34").Value = 1470.8223
$ws.Range("L34").Value = 33337456
$ws.Range("M34").Value = -1268.8223
$ws.Range("N34").Value = -33337860
$ws.Range("H58").Value = 1193513.4
$ws.Range("I58").Value = 2208
$ws.Range("J58").Value = 3576124.2
$ws.Range("K58").Value = 2208
$ws.Range("L58").Value = 3576124.2
$ws.Range("M58").Value = -2005
$ws.Range("N58").Value = -3576530.2
$ws.Range("H107").Value = 558.439
$ws.Range("I107").Value = 243.95
$ws.Range("J107").Value = 857.9524
$ws.Range("K107").Value = 243.95
$ws.Range("L107").Value = 857.9524
$ws.Range("M107").Value = 1676.05
$ws.Range("N107").Value = -4697.9524
$ws.Range("H113").Value = 4911.1
$ws.Range("I113").Value = 5302.2
$ws.Range("J113").Value = 4520
$ws.Range("K113").Value = 5302.2
$ws.Range("L113").Value = 4520
$ws.Range("M113").Value = -3132.2
$ws.Range("N113").Value = -8860
$ws.Range("H122").Value = 43479556
$ws.Range("I122").Value = 76924020
$ws.Range("J122").Value = 1755
$ws.Range("K122").Value = 230772060
$ws.Range("L122").Value = 5265
$ws.Range("M122").Value = -230769610
$ws.Range("N122").Value = -10165
$ws.Range("H132").Value = 2858.6667
$ws.Range("I132").Value = 1749.5385
$ws.Range("J132").Value = 3706.8235
$ws.Range("K132").Value = 5248.6155
$ws.Range("L132").Value = 11120.4705
$ws.Range("M132").Value = -2718.6155
$ws.Range("N132").Value = -16180.4705
$ws.Range("H136").Value = 1193513.4
$ws.Range("I136").Value = 2208
$ws.Range("J136").Value = 3576124.2
$ws.Range("K136").Value = 6624
$ws.Range("L136").Value = 10728372.6
$ws.Range("M136").Value = -4074
$ws.Range("N136").Value = -10733472.6
$ws.Range("H138").Value = 45184
$ws.Range("J138").Value = 45184
$ws.Range("L138").Value = 45184
$ws.Range("N138").Value = -55464
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 686.96
$ws.Range("I5").Value = 436.31818
$ws.Range("J5").Value = 2525
$ws.Range("K5").Value = 1308.95454
$ws.Range("L5").Value = 7575
$ws.Range("M5").Value = -1196.95454
$ws.Range("N5").Value = -7799
$ws.Range("H117").Value = 10205956
$ws.Range("I117").Value = 413.85715
$ws.Range("J117").Value = 20411498
$ws.Range("K117").Value = 1241.57145
$ws.Range("L117").Value = 61234494
$ws.Range("M117").Value = 2200.42855
$ws.Range("N117").Value = -61241378
$ws.Range("H122").Value = 3111.6897
$ws.Range("I122").Value = 418
$ws.Range("J122").Value = 3672.875
$ws.Range("K122").Value = 3762
$ws.Range("L122").Value = 33055.875
$ws.Range("M122").Value = -1312
$ws.Range("N122").Value = -37955.875
$ws.Range("H132").Value = 3645.4
$ws.Range("I132").Value = 1725
$ws.Range("K132").Value = 15525
$ws.Range("M132").Value = -12995
$ws.Range("H135").Value = 686.96
$ws.Range("I135").Value = 436.31818
$ws.Range("J135").Value = 2525
$ws.Range("K135").Value = 3926.86362
$ws.Range("L135").Value = 22725
$ws.Range("M135").Value = -1391.86362
$ws.Range("N135").Value = -27795
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 790.5625
$ws.Range("I107").Value = 782.8889
$ws.Range("J107").Value = 800.4286
$ws.Range("K107").Value = 782.8889
$ws.Range("L107").Value = 800.4286
$ws.Range("M107").Value = 1137.1111
$ws.Range("N107").Value = -4640.4286
$ws.Range("H113").Value = 3853.25
$ws.Range("I113").Value = 3750
$ws.Range("J113").Value = 3956.5
$ws.Range("K113").Value = 3750
$ws.Range("L113").Value = 3956.5
$ws.Range("M113").Value = -1580
$ws.Range("N113").Value = -8296.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
